# Updated cryptos list with GitHub Actions run.
# Column D holds price strings that look numeric (e.g. "246.48"); force the
# whole price column to Text format first so the new values are written back
# as literal strings instead of being auto-coerced to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.329.16"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.882.39"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "246.48"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Value = "43.45"
$ws.Range("E8").Value = "  +5.34%  "
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "53.29"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "13.51"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").Value = "2.154.97"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "1.902.18"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "35.336.63"
$ws.Range("D19").Value = "73.80"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "244.65"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("E24").Value = "  +8.37%  "
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "164.82"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "18.30"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "0.0595"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").Value = "4.20"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -8.18%  "
$ws.Range("D37").Value = "0.854"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("E39").Value = "  +10.24%  "
$ws.Range("D40").Value = "17.55"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").Value = "97.08"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "1.309.38"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "0.0801"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").Value = "11.78"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "6.33"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "42.38"
$ws.Range("E51").Value = "  +0.15%  "
